$wb = $excel.ActiveWorkbook

$commitHash = "97b21b0da07e90d7439ad8308bb0d882cd3f6104"
$fileGuid1 = "6b7eca85-f4a3-4bcf-a4ff-722720788659"
$fileGuid2 = "81015c87-6f0f-49eb-bec6-73ae877e7b2e"
$zhHash1   = "c0f883fac2efdaf8506e21f31255ff67ec8794bd"
$zhHash2   = "1fc56a5d2d6d58e5cf743b2764e018930a77500a"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Cells.Item(4,1).Value = "$fileGuid1.md"
$wsOverview.Cells.Item(4,2).Value = "e2e\$fileGuid1.md"
$wsOverview.Cells.Item(4,3).Value = ".md"
$wsOverview.Cells.Item(4,5).Value = "Ready for handoff"
$wsOverview.Cells.Item(4,6).Value = "Ready for handoff"
$wsOverview.Cells.Item(4,7).Value = "2016-10-26 07:13:31"

$wsOverview.Cells.Item(5,1).Value = "$fileGuid2.md"
$wsOverview.Cells.Item(5,2).Value = "e2e\$fileGuid2.md"
$wsOverview.Cells.Item(5,3).Value = ".md"
$wsOverview.Cells.Item(5,5).Value = "Ready for handoff"
$wsOverview.Cells.Item(5,6).Value = "Ready for handoff"
$wsOverview.Cells.Item(5,7).Value = "2016-10-26 07:13:31"

$wsOverview.Hyperlinks.Add($wsOverview.Cells.Item(4,2), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commitHash/e2e/$fileGuid1.md", [Type]::Missing, [Type]::Missing, "e2e\$fileGuid1.md")
$wsOverview.Hyperlinks.Add($wsOverview.Cells.Item(5,2), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commitHash/e2e/$fileGuid2.md", [Type]::Missing, [Type]::Missing, "e2e\$fileGuid2.md")

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G5"))

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Cells.Item(4,1).Value = "$fileGuid1.md"
$wsZh.Cells.Item(4,2).Value = ".md"
$wsZh.Cells.Item(4,3).Value = "Ready for handoff"
$wsZh.Cells.Item(4,4).Value = "e2e"
$wsZh.Cells.Item(4,5).Value = "ht"
$wsZh.Cells.Item(4,6).Value = "'False"
$wsZh.Cells.Item(4,7).Value = "$fileGuid1.$zhHash1.zh-cn.xlf"
$wsZh.Cells.Item(4,8).Value = "2016-10-26 07:13:20"
$wsZh.Cells.Item(4,11).Value = "0001-01-01 00:00:00"
$wsZh.Cells.Item(4,13).Value = "'True"
$wsZh.Cells.Item(4,15).Value = "'False"

$wsZh.Cells.Item(5,1).Value = "$fileGuid2.md"
$wsZh.Cells.Item(5,2).Value = ".md"
$wsZh.Cells.Item(5,3).Value = "Ready for handoff"
$wsZh.Cells.Item(5,4).Value = "e2e"
$wsZh.Cells.Item(5,5).Value = "ht"
$wsZh.Cells.Item(5,6).Value = "'False"
$wsZh.Cells.Item(5,7).Value = "$fileGuid2.$zhHash2.zh-cn.xlf"
$wsZh.Cells.Item(5,8).Value = "2016-10-26 07:13:20"
$wsZh.Cells.Item(5,11).Value = "0001-01-01 00:00:00"
$wsZh.Cells.Item(5,13).Value = "'True"
$wsZh.Cells.Item(5,15).Value = "'False"

$wsZh.Hyperlinks.Add($wsZh.Cells.Item(4,1), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commitHash/e2e/$fileGuid1.md", [Type]::Missing, [Type]::Missing, "$fileGuid1.md")
$wsZh.Hyperlinks.Add($wsZh.Cells.Item(5,1), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commitHash/e2e/$fileGuid2.md", [Type]::Missing, [Type]::Missing, "$fileGuid2.md")

$loZh = $wsZh.ListObjects.Item(1)
$loZh.Resize($wsZh.Range("A1:P5"))

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Cells.Item(4,1).Value = "$fileGuid1.md"
$wsDe.Cells.Item(4,2).Value = ".md"
$wsDe.Cells.Item(4,3).Value = "Ready for handoff"
$wsDe.Cells.Item(4,4).Value = "e2e"
$wsDe.Cells.Item(4,5).Value = "ht"
$wsDe.Cells.Item(4,6).Value = "False"
$wsDe.Cells.Item(4,7).Value = "$fileGuid1.$zhHash1.de-de.xlf"
$wsDe.Cells.Item(4,8).Value = "2016-10-26 07:13:31"
$wsDe.Cells.Item(4,11).Value = "0001-01-01 00:00:00"
$wsDe.Cells.Item(4,13).Value = "True"
$wsDe.Cells.Item(4,15).Value = "False"

$wsDe.Cells.Item(5,1).Value = "$fileGuid2.md"
$wsDe.Cells.Item(5,2).Value = ".md"
$wsDe.Cells.Item(5,3).Value = "Ready for handoff"
$wsDe.Cells.Item(5,4).Value = "e2e"
$wsDe.Cells.Item(5,5).Value = "ht"
$wsDe.Cells.Item(5,6).Value = "False"
$wsDe.Cells.Item(5,7).Value = "$fileGuid2.$zhHash2.de-de.xlf"
$wsDe.Cells.Item(5,8).Value = "2016-10-26 07:13:31"
$wsDe.Cells.Item(5,11).Value = "0001-01-01 00:00:00"
$wsDe.Cells.Item(5,13).Value = "True"
$wsDe.Cells.Item(5,15).Value = "False"

$wsDe.Hyperlinks.Add($wsDe.Cells.Item(4,1), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commitHash/e2e/$fileGuid1.md", [Type]::Missing, [Type]::Missing, "$fileGuid1.md")
$wsDe.Hyperlinks.Add($wsDe.Cells.Item(5,1), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commitHash/e2e/$fileGuid2.md", [Type]::Missing, [Type]::Missing, "$fileGuid2.md")

$loDe = $wsDe.ListObjects.Item(1)
$loDe.Resize($wsDe.Range("A1:P5"))

# ---------------------------------------------------------------------------
# Column width tweaks (status/language columns widened to fit "Ready for handoff")
# ---------------------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 17.2
$wsOverview.Columns.Item(6).ColumnWidth = 17.2
$wsZh.Columns.Item(3).ColumnWidth = 17.2
$wsDe.Columns.Item(3).ColumnWidth = 17.2
